$wb = $excel.ActiveWorkbook

# Turn on "Include Emissions from Imported Electricity" boolean lever (B2: 0 -> 1)
$ws = $wb.Worksheets.Item("BIEfIE")
$ws.Range("B2").Value = 1

# Make the BIEfIE sheet the active/selected sheet, with B3 selected
$ws.Activate()
$ws.Range("B3").Select()
